$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

# Calibration tweak: lower the "share of costs that must be covered"
# values in B2:B25 from 1 (100%) to 0.95 (95%)
$ws.Range("B2:B25").Value = 0.95

# Reflect the new selection left in the sheet (B2:B25, active cell B2)
# without changing which sheet tab is active in the workbook.
$activeSheetName = $wb.ActiveSheet.Name
$ws.Select()
$ws.Range("B2:B25").Select()
$wb.Worksheets.Item($activeSheetName).Select()
